$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Elimina antiguos EC" - the old table had 6 detail rows (periods 2506, 2506, 2507, 2507, 2508, 2508).
# Carry the "bottom of table" border formatting (currently on row 21) up onto row 17, which will
# become the new last row of the table once the old period rows (18-21) are removed.
$ws.Range("B21:J21").Copy() | Out-Null
$ws.Range("B17:J17").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Remove the old rows for periods 2507 and 2508 (rows 18-21)
$ws.Range("18:21").Delete()

# "agrega nuevos" - the remaining two detail rows (now rows 16-17) report the new period 2509
$ws.Range("E16").Value = "2509"
$ws.Range("E16").HorizontalAlignment = -4108  # xlHAlignCenter
$ws.Range("E17").Value = "2509"
$ws.Range("E17").HorizontalAlignment = -4108  # xlHAlignCenter

# "modifica antigua BD" - update the summary figures
$ws.Range("E11").Value = 113880
$ws.Range("F13").Value = 1
